$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename TestBean -> JavaBean in the two cells that reference it.
$ws.Range("F5").Value = "Data JavaBean beans2"
$ws.Range("B3").Value = "Method String print(JavaBean bean)"

# Move the active selection to B4 (as in the saved workbook state).
$ws.Range("B4").Select() | Out-Null

# Normalize the theme naming (Russian-locale "Стандартная" -> "Office") that
# accompanies the resave.
$wb.Theme.ThemeFontScheme.Name = "Office"
